# Weekly update of the "Fruta - Tuna" price table:
# a new weekly record (Especial/Primera, date 44636) is inserted at the top
# of the historical block (rows 38-58), pushing the older weekly groups
# down by two rows; the two rows that fall off the bottom of the original
# range are re-appended at the end (rows 59-60), growing the sheet from
# A1:T58 to A1:T60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are constant across every data row of this subset.
$constCols = @{
    A = 2
    B = "Comercializadora del Agro de Limarí"
    C = "Coquimbo"
    E = 4
    F = "Fruta"
    G = 100107
    H = "Otros"
    I = 100107011
    J = "Tuna"
    K = "Sin especificar"
    Q = "$/caja 18 kilos"
    R = "Provincia de Limarí"
    T = 18
}

# Final state (rows 38-60) for the columns that actually vary: D (fecha),
# L (calidad), M (volumen), N (precio min), O (precio max), P (precio
# promedio ponderado), S (precio $/kg).
$rows = @(
    @{Row=38; D=44636; L="Especial"; M=300; N=11000; O=12000; P=11500; S=639},
    @{Row=39; D=44636; L="Primera"; M=300; N=9000; O=10000; P=9500; S=528},
    @{Row=40; D=44280; L="Especial"; M=240; N=12500; O=13000; P=12750; S=708},
    @{Row=41; D=44280; L="Primera"; M=240; N=10500; O=11000; P=10750; S=597},
    @{Row=42; D=44280; L="Segunda"; M=300; N=8500; O=9000; P=8750; S=486},
    @{Row=43; D=44294; L="Especial"; M=360; N=12500; O=13000; P=12750; S=708},
    @{Row=44; D=44294; L="Primera"; M=240; N=10500; O=11000; P=10750; S=597},
    @{Row=45; D=44294; L="Segunda"; M=240; N=8500; O=9000; P=8750; S=486},
    @{Row=46; D=44279; L="Especial"; M=200; N=12500; O=13000; P=12750; S=708},
    @{Row=47; D=44279; L="Primera"; M=240; N=10500; O=11000; P=10750; S=597},
    @{Row=48; D=44279; L="Segunda"; M=240; N=8500; O=9000; P=8750; S=486},
    @{Row=49; D=44385; L="Especial"; M=120; N=14000; O=14500; P=14250; S=792},
    @{Row=50; D=44385; L="Primera"; M=300; N=11000; O=11500; P=11250; S=625},
    @{Row=51; D=44385; L="Segunda"; M=240; N=8000; O=8500; P=8250; S=458},
    @{Row=52; D=44385; L="Tercera"; M=120; N=5000; O=5500; P=5250; S=292},
    @{Row=53; D=44272; L="Especial"; M=160; N=12500; O=13000; P=12750; S=708},
    @{Row=54; D=44272; L="Primera"; M=300; N=10500; O=11000; P=10750; S=597},
    @{Row=55; D=44272; L="Segunda"; M=240; N=8500; O=9000; P=8750; S=486},
    @{Row=56; D=44615; L="Especial"; M=200; N=14000; O=15000; P=14500; S=806},
    @{Row=57; D=44615; L="Primera"; M=400; N=12000; O=13000; P=12500; S=694},
    @{Row=58; D=44335; L="Especial"; M=240; N=19500; O=20000; P=19750; S=1097},
    @{Row=59; D=44335; L="Primera"; M=200; N=17500; O=18000; P=17750; S=986},
    @{Row=60; D=44335; L="Segunda"; M=160; N=12500; O=13000; P=12750; S=708}
)

foreach ($item in $rows) {
    $r = $item.Row

    # Constant template columns — only need to be (re)written for the two
    # brand-new rows (59, 60); harmless no-ops for the pre-existing rows.
    foreach ($col in $constCols.Keys) {
        $ws.Range("$col$r").Value = $constCols[$col]
    }

    $ws.Range("D$r").Value = $item.D
    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
    $ws.Range("O$r").Value = $item.O
    $ws.Range("P$r").Value = $item.P
    $ws.Range("S$r").Value = $item.S
}
